# Insert a new weekly price record as row 169, pushing the existing
# rows 169-272 down to 170-273 (the sheet keeps collecting new weekly
# observations, so the most recent record lands at the top of the
# block and everything else shifts down by one row).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("169:169").Insert()

$ws.Cells.Item(169, 1).Value = 8
$ws.Cells.Item(169, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(169, 3).Value = "Coquimbo"
$ws.Cells.Item(169, 4).Value = 44582
$ws.Cells.Item(169, 5).Value = 4
$ws.Cells.Item(169, 6).Value = 100114013
$ws.Cells.Item(169, 7).Value = "Zanahoria"
$ws.Cells.Item(169, 8).Value = "Sin especificar"
$ws.Cells.Item(169, 9).Value = "Primera"
$ws.Cells.Item(169, 10).Value = 780
$ws.Cells.Item(169, 11).Value = 6000
$ws.Cells.Item(169, 12).Value = 6500
$ws.Cells.Item(169, 13).Value = 6250
$ws.Cells.Item(169, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(169, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(169, 16).Value = 312
$ws.Cells.Item(169, 17).Value = 20
$ws.Cells.Item(169, 18).Value = "Hortaliza"
